$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 533.8889
$ws.Range("I33").Value = 535.2692
$ws.Range("K33").Value = 535.2692
$ws.Range("M33").Value = -306.2692
$ws.Range("H40").Value = 1969060
$ws.Range("I40").Value = 15437.375
$ws.Range("J40").Value = 3705613.5
$ws.Range("K40").Value = 15437.375
$ws.Range("L40").Value = 3705613.5
$ws.Range("M40").Value = -15262.375
$ws.Range("N40").Value = -3705963.5
$ws.Range("H58").Value = 3834.111
$ws.Range("I58").Value = 998
$ws.Range("J58").Value = 7379.25
$ws.Range("K58").Value = 2994
$ws.Range("L58").Value = 22137.75
$ws.Range("M58").Value = -2844
$ws.Range("N58").Value = -22437.75
$ws.Range("H61").Value = 41668680
$ws.Range("I61").Value = 41668680
$ws.Range("K61").Value = 125006040
$ws.Range("M61").Value = -125005868
$ws.Range("H62").Value = 55572240
$ws.Range("I62").Value = 111113020
$ws.Range("J62").Value = 31452.889
$ws.Range("K62").Value = 111113020
$ws.Range("L62").Value = 31452.889
$ws.Range("M62").Value = -111112396
$ws.Range("N62").Value = -32700.889
$ws.Range("H65").Value = 55572240
$ws.Range("I65").Value = 111113020
$ws.Range("J65").Value = 31452.889
$ws.Range("K65").Value = 555565100
$ws.Range("L65").Value = 157264.445
$ws.Range("M65").Value = -555561980
$ws.Range("N65").Value = -163504.445
$ws.Range("H76").Value = 17075
$ws.Range("J76").Value = 15000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15630
$ws.Range("H79").Value = 17075
$ws.Range("J79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17184
$ws.Range("H106").Value = 1914.0526
$ws.Range("I106").Value = 1914.0526
$ws.Range("K106").Value = 1914.0526
$ws.Range("M106").Value = -1283.0526
$ws.Range("H132").Value = 1910.4546
$ws.Range("I132").Value = 1813.9375
$ws.Range("K132").Value = 5441.8125
$ws.Range("M132").Value = -2911.8125
$ws.Range("H137").Value = 2515.963
$ws.Range("I137").Value = 2561.6667
$ws.Range("J137").Value = 2458.8333
$ws.Range("K137").Value = 7685.000100000001
$ws.Range("L137").Value = 7376.499899999999
$ws.Range("M137").Value = -5135.000100000001
$ws.Range("N137").Value = -12476.4999
$ws.Range("H138").Value = 1544384.2
$ws.Range("I138").Value = 4137.75
$ws.Range("J138").Value = 1893119.2
$ws.Range("K138").Value = 12413.25
$ws.Range("L138").Value = 5679357.6
$ws.Range("M138").Value = -7273.25
$ws.Range("N138").Value = -5689637.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3344.0984
$ws.Range("I32").Value = 3344.0984
$ws.Range("K32").Value = 3344.0984
$ws.Range("M32").Value = -3057.0984
$ws.Range("H61").Value = 5908.9463
$ws.Range("I61").Value = 3461.2683
$ws.Range("K61").Value = 3461.2683
$ws.Range("M61").Value = -3249.2683
$ws.Range("H63").Value = 2500
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814
$ws.Range("H66").Value = 2500
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068
$ws.Range("H80").Value = 30050
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 30050
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H122").Value = 6111.4707
$ws.Range("I122").Value = 7838.857
$ws.Range("K122").Value = 23516.571
$ws.Range("M122").Value = -21066.571
$ws.Range("H136").Value = 5908.9463
$ws.Range("I136").Value = 3461.2683
$ws.Range("K136").Value = 10383.8049
$ws.Range("M136").Value = -7833.804900000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4660.143
$ws.Range("I134").Value = 1654.9048
$ws.Range("K134").Value = 4964.7144
$ws.Range("M134").Value = -2429.7144

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7216.875
$ws.Range("I99").Value = 7632.4443
$ws.Range("J99").Value = 6967.533
$ws.Range("K99").Value = 7632.4443
$ws.Range("L99").Value = 6967.533
$ws.Range("M99").Value = -6134.4443
$ws.Range("N99").Value = -9963.532999999999
$ws.Range("H126").Value = 7216.875
$ws.Range("I126").Value = 7632.4443
$ws.Range("J126").Value = 6967.533
$ws.Range("K126").Value = 22897.3329
$ws.Range("L126").Value = 20902.599
$ws.Range("M126").Value = -20427.3329
$ws.Range("N126").Value = -25842.599
$ws.Range("H132").Value = 3754.353
$ws.Range("I132").Value = 2245.5833
$ws.Range("J132").Value = 7375.4
$ws.Range("K132").Value = 6736.749899999999
$ws.Range("L132").Value = 22126.2
$ws.Range("M132").Value = -4206.749899999999
$ws.Range("N132").Value = -27186.2
$ws.Range("H141").Value = 70190.25
$ws.Range("J141").Value = 70190.25
$ws.Range("L141").Value = 70190.25
$ws.Range("N141").Value = -80550.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 1217.625
$ws.Range("J61").Value = 1924.5
$ws.Range("L61").Value = 5773.5
$ws.Range("N61").Value = -6203.5
$ws.Range("H107").Value = 22223488
$ws.Range("J107").Value = 28572842
$ws.Range("L107").Value = 85718526
$ws.Range("N107").Value = -85722366
$ws.Range("H113").Value = 2551.1724
$ws.Range("I113").Value = 1040
$ws.Range("J113").Value = 3126.8572
$ws.Range("K113").Value = 3120
$ws.Range("L113").Value = 9380.571599999999
$ws.Range("M113").Value = -950
$ws.Range("N113").Value = -13720.5716
$ws.Range("H122").Value = 708354.0600000001
$ws.Range("I122").Value = 1768992.2
$ws.Range("J122").Value = 1262
$ws.Range("K122").Value = 15920929.8
$ws.Range("L122").Value = 11358
$ws.Range("M122").Value = -15918479.8
$ws.Range("N122").Value = -16258
$ws.Range("H132").Value = 9536.652
$ws.Range("J132").Value = 11119.6
$ws.Range("L132").Value = 100076.4
$ws.Range("N132").Value = -105136.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H11").Value = 1100000
$ws.Range("I11").Value = 1000000
$ws.Range("J11").Value = 1200000
$ws.Range("K11").Value = 1000000
$ws.Range("L11").Value = 1200000
$ws.Range("M11").Value = -999861
$ws.Range("N11").Value = -1200278
$ws.Range("H19").Value = 3002.5
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = 6000
$ws.Range("K19").Value = 5
$ws.Range("L19").Value = 6000
$ws.Range("M19").Value = 283
$ws.Range("N19").Value = -6576
$ws.Range("H80").Value = 3059.923
$ws.Range("J80").Value = 3081
$ws.Range("L80").Value = 3081
$ws.Range("N80").Value = -5077
$ws.Range("H83").Value = 3059.923
$ws.Range("J83").Value = 3081
$ws.Range("L83").Value = 15405
$ws.Range("N83").Value = -25389
$ws.Range("H122").Value = 105499.4
$ws.Range("I122").Value = 500999.5
$ws.Range("J122").Value = 6624.375
$ws.Range("K122").Value = 1502998.5
$ws.Range("L122").Value = 19873.125
$ws.Range("M122").Value = -1500548.5
$ws.Range("N122").Value = -24773.125
$ws.Range("H132").Value = 7829.7
$ws.Range("I132").Value = 6040
$ws.Range("K132").Value = 18120
$ws.Range("M132").Value = -15590

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5525.9414
$ws.Range("I40").Value = 5209.5
$ws.Range("K40").Value = 5209.5
$ws.Range("M40").Value = -5073.5
$ws.Range("H46").Value = 11498921
$ws.Range("J46").Value = 7001.5
$ws.Range("L46").Value = 7001.5
$ws.Range("N46").Value = -7377.5
$ws.Range("H55").Value = 66667190
$ws.Range("I55").Value = 111111400
$ws.Range("K55").Value = 111111400
$ws.Range("M55").Value = -111111227
$ws.Range("H122").Value = 4978.5137
$ws.Range("I122").Value = 4355.3184
$ws.Range("J122").Value = 5892.533
$ws.Range("K122").Value = 13065.9552
$ws.Range("L122").Value = 17677.599
$ws.Range("M122").Value = -10615.9552
$ws.Range("N122").Value = -22577.599

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 75833.164
$ws.Range("J46").Value = 75833.164
$ws.Range("L46").Value = 75833.164
$ws.Range("N46").Value = -76295.164
$ws.Range("H122").Value = 2635.4666
$ws.Range("I122").Value = 2063.3416
$ws.Range("K122").Value = 6190.024800000001
$ws.Range("M122").Value = -3740.024800000001
$ws.Range("H132").Value = 11123764
$ws.Range("I132").Value = 13521062
$ws.Range("J132").Value = 36260.875
$ws.Range("K132").Value = 40563186
$ws.Range("L132").Value = 108782.625
$ws.Range("M132").Value = -40560656
$ws.Range("N132").Value = -113842.625
$ws.Range("H134").Value = 75833.164
$ws.Range("J134").Value = 75833.164
$ws.Range("L134").Value = 227499.492
$ws.Range("N134").Value = -232569.492
